$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 update: refresh last_action_date, bump reactions/replies counts,
# and append new message ids to the tracked lists.
$ws.Range("E21").Value = "2026-02-13T07:27:50.964363+00:00"
$ws.Range("H21").Value = 2
$ws.Range("I21").Value = 2
$ws.Range("L21").Value = "[19, 32]"
$ws.Range("M21").Value = "[27, 14]"
